# Update Price (D) and Volume(1h) (E) columns per the latest crypto data refresh.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'28.459.50"
$ws.Range("E2").Value = "  -0.02%  "
$ws.Range("D3").Value = "'1.824.42"
$ws.Range("E3").Value = "  -0.23%  "
$ws.Range("D4").Value = "'1.004"
$ws.Range("E4").Value = "  +0.23%  "
$ws.Range("D5").Value = "'316.95"
$ws.Range("E5").Value = "  +0.54%  "
$ws.Range("E6").Value = "  +0.18%  "
$ws.Range("D7").Value = "'0.5168"
$ws.Range("E7").Value = "  +0.79%  "
$ws.Range("D8").Value = "'0.3858"
$ws.Range("E8").Value = "  -1.24%  "
$ws.Range("E9").Value = "  +8.65%  "
$ws.Range("D10").Value = "'1.123"
$ws.Range("E10").Value = "  +1.32%  "
$ws.Range("D11").Value = "'41.89"
$ws.Range("E11").Value = "  -0.01%  "
$ws.Range("D12").Value = "'6.377"
$ws.Range("E12").Value = "  +1.33%  "
$ws.Range("E13").Value = "  -0.01%  "
$ws.Range("E14").Value = "  +0.16%  "
$ws.Range("D15").Value = "'7.473"
$ws.Range("E15").Value = "  -1.02%  "
$ws.Range("D16").Value = "'1.823.85"
$ws.Range("E16").Value = "  -0.03%  "
$ws.Range("D17").Value = "'94.22"
$ws.Range("E17").Value = "  +0.87%  "
$ws.Range("D18").Value = "'0.00001122"
$ws.Range("E18").Value = "  +3.69%  "
$ws.Range("D19").Value = "'0.06635"
$ws.Range("E19").Value = "  -0.41%  "
$ws.Range("E20").Value = "  +0.42%  "
$ws.Range("E21").Value = "  +0.20%  "
$ws.Range("D22").Value = "'6.051"
$ws.Range("E22").Value = "  -2.15%  "
$ws.Range("D23").Value = "'28.489.78"
$ws.Range("E23").Value = "  -0.01%  "
$ws.Range("D24").Value = "'11.49"
$ws.Range("E24").Value = "  +2.91%  "
$ws.Range("D25").Value = "'2.243"
$ws.Range("D26").Value = "'21.11"
$ws.Range("E26").Value = "  +2.22%  "
$ws.Range("D27").Value = "'159.35"
$ws.Range("E27").Value = "  +1.62%  "
$ws.Range("D28").Value = "'2.034.14"
$ws.Range("E28").Value = "  -0.02%  "
$ws.Range("D29").Value = "'2.406"
$ws.Range("E29").Value = "  +0.42%  "
$ws.Range("D30").Value = "'126.02"
$ws.Range("E30").Value = "  +0.72%  "
$ws.Range("D31").Value = "'0.1110"
$ws.Range("E31").Value = "  +2.26%  "
$ws.Range("D32").Value = "'1.092"
$ws.Range("E32").Value = "  -2.94%  "
$ws.Range("D33").Value = "'5.730"
$ws.Range("E33").Value = "  +0.99%  "
$ws.Range("D34").Value = "'0.07522"
$ws.Range("E34").Value = "  +7.10%  "
$ws.Range("D35").Value = "'3.682"
$ws.Range("E35").Value = "  +0.49%  "
$ws.Range("D36").Value = "'0.2225"
$ws.Range("E36").Value = "  -0.17%  "
$ws.Range("D37").Value = "'0.02361"
$ws.Range("E37").Value = "  +1.61%  "
$ws.Range("E38").Value = "  +7.67%  "
$ws.Range("D39").Value = "'5.260"
$ws.Range("E39").Value = "  +2.36%  "
$ws.Range("D40").Value = "'8.757"
$ws.Range("E40").Value = "  -2.71%  "
$ws.Range("D41").Value = "'0.6389"
$ws.Range("E41").Value = "  +1.80%  "
$ws.Range("D42").Value = "'1.186"
$ws.Range("E42").Value = "  +0.29%  "
$ws.Range("D43").Value = "'1.395"
$ws.Range("E43").Value = "  -0.26%  "
$ws.Range("D44").Value = "'0.6194"
$ws.Range("E44").Value = "  +4.92%  "
$ws.Range("D45").Value = "'13.56"
$ws.Range("E45").Value = "  +0.62%  "
$ws.Range("D46").Value = "'3.798"
$ws.Range("E46").Value = "  +2.25%  "
$ws.Range("D47").Value = "'127.60"
$ws.Range("E47").Value = "  +2.69%  "
$ws.Range("E48").Value = "  +1.08%  "
$ws.Range("D49").Value = "'1.205"
$ws.Range("E49").Value = "  +0.48%  "
$ws.Range("D50").Value = "'0.06963"
$ws.Range("E50").Value = "  +0.57%  "
$ws.Range("D51").Value = "'1.082"
$ws.Range("E51").Value = "  +1.44%  "
